$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Re-use the existing "big decorative font" styles (hiragana tfb / bold
#    Calibri) for the row labels in column A, since those decorative fonts
#    used to be applied to the letter/number answer cells which now hold
#    picture placeholders instead of rendered text ("pictures instead of
#    letters").
# ---------------------------------------------------------------------------

# A2 ("Most Probable") takes on the big "hiragana tfb" font previously used
# for the rendered-letters cell (B2).
$ws.Range("B2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

# A3 ("Least Probable") takes on the big bold Calibri font previously used
# for the rendered-numbers cell (C2/C3).
$ws.Range("C2").Copy()
$ws.Range("A3").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# The former letter/number cells go back to the plain default look, since
# they now just hold a picture filename placeholder.
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Style = "Normal"
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Style = "Normal"

# ---------------------------------------------------------------------------
# 2) Update the cell text. B/C columns on rows 2 and 3 now reference picture
#    file names instead of literal letters/numbers.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "question"

$ws.Range("A2").Value = "Most Probable"
$ws.Range("B2").Value = "Hselect1.jpg"
$ws.Range("C2").Value = "Nselect1.jpg"

$ws.Range("A3").Value = "Least Probable"
$ws.Range("B3").Value = "Hselect1.jpg"
$ws.Range("C3").Value = "Nselect1.jpg"

# ---------------------------------------------------------------------------
# 3) Rows shrink now that the big decorative letter/number glyphs are gone
#    from columns B/C.
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 37.5
$ws.Rows.Item(3).RowHeight = 37.5

# ---------------------------------------------------------------------------
# 4) Misc view / print settings.
# ---------------------------------------------------------------------------
[void]$ws.Range("B7").Select()

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
